$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Pravin"
$ws.Range("B2").Value = "Paul"
$ws.Range("C2").Value = "Raj"

$ws.Range("C3").Select() | Out-Null
